$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new data row (row 2) with the collection's metadata.
$ws.Range("A2").Value = "MCH117-1"
$ws.Range("C2").Value = "PAPERS"
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: 21B | GRAP COUNT NUMER: NONE"

# Match the formatting used for the rest of the new row: 10pt Calibri,
# theme text color (matches the font used elsewhere in the sheet).
$fmtRange = $ws.Range("A2:H2")
$fmtRange.Font.Name = "Calibri"
$fmtRange.Font.Size = 10
$fmtRange.Font.ThemeColor = 1

# Restore the selection/frozen-pane view state on the new row.
$ws.Range("A2:I2").Select()
$win = $excel.ActiveWindow
$win.FreezePanes = $true
